$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for two new rows above the old "RCP85" model-name row -------
# The old row 2 (RCP85 + model run names) needs to end up at row 5, with a
# new "In Paper" (Y/N) row at row 2, a numbering row at row 3, and row 4
# left blank. Inserting 3 rows at row 2 shifts the old row 2 down to row 5.
$ws.Rows("2:4").Insert()

# --- Row 6: "RCP8.5 In paper" indicator row --------------------------------
# Fill the Y's first (so the shared string "Y" is created before the longer
# label), then the row label itself.
$ws.Range("B6").Value2 = "Y"
$ws.Range("C6").Value2 = "Y"
$ws.Range("D6").Value2 = "Y"
$ws.Range("E6").Value2 = "Y"
$ws.Range("F6").Value2 = "Y"
$ws.Range("G6").Value2 = "Y"
$ws.Range("H6").Value2 = "Y"
$ws.Range("I6").Value2 = "Y"
$ws.Range("J6").Value2 = "Y"
$ws.Range("K6").Value2 = "Y"
$ws.Range("L6").Value2 = "Y"
$ws.Range("M6").Value2 = "Y"
$ws.Range("N6").Value2 = "Y"
$ws.Range("O6").Value2 = "Y"
$ws.Range("P6").Value2 = "Y"
$ws.Range("Q6").Value2 = "Y"
$ws.Range("R6").Value2 = "Y"
$ws.Range("S6").Value2 = "Y"
$ws.Range("T6").Value2 = "Y"
$ws.Range("U6").Value2 = "Y"
$ws.Range("V6").Value2 = "Y"
$ws.Range("A6").Value2 = "RCP8.5 In paper"

# --- Row 2: "In Paper" (Y/N) header row -------------------------------------
# Again fill the Y's before the N's and before the row label, matching the
# original authoring order (Y, then RCP8.5 In paper, then N, then In Paper).
$ws.Range("B2").Value2 = "Y"
$ws.Range("C2").Value2 = "Y"
$ws.Range("D2").Value2 = "Y"
$ws.Range("E2").Value2 = "Y"
$ws.Range("F2").Value2 = "Y"
$ws.Range("G2").Value2 = "Y"
$ws.Range("H2").Value2 = "Y"
$ws.Range("I2").Value2 = "Y"
$ws.Range("J2").Value2 = "Y"
$ws.Range("K2").Value2 = "Y"
$ws.Range("L2").Value2 = "Y"
$ws.Range("M2").Value2 = "Y"
$ws.Range("P2").Value2 = "Y"
$ws.Range("Q2").Value2 = "Y"
$ws.Range("R2").Value2 = "Y"
$ws.Range("S2").Value2 = "Y"
$ws.Range("T2").Value2 = "Y"
$ws.Range("U2").Value2 = "Y"
$ws.Range("V2").Value2 = "Y"
$ws.Range("X2").Value2 = "Y"
$ws.Range("Y2").Value2 = "Y"
$ws.Range("N2").Value2 = "N"
$ws.Range("O2").Value2 = "N"
$ws.Range("W2").Value2 = "N"
$ws.Range("A2").Value2 = "In Paper"

# --- Row 3: plain column numbering 1-24 in B3:Y3 ----------------------------
for ($i = 2; $i -le 25; $i++) {
    $ws.Cells.Item(3, $i).Value2 = $i - 1
}

# --- Selection / scroll position matching the saved view -------------------
$ws.Range("V3").Select()

# --- Column widths -----------------------------------------------------------
# A and C were manually widened; D:Y were set via "AutoFit Column Width".
# ColumnWidth is quantized to 1/6-character increments by this host, so the
# values below are the closest representable widths to the authored ones.
$ws.Columns("A").ColumnWidth = 13.8333333333333
$ws.Columns("C").ColumnWidth = 14.5
$ws.Columns("D").ColumnWidth = 13.3333333333333
$ws.Columns("E:G").ColumnWidth = 13.1666666666667
$ws.Columns("H").ColumnWidth = 14.1666666666667
$ws.Columns("I").ColumnWidth = 12
$ws.Columns("J").ColumnWidth = 12.3333333333333
$ws.Columns("K").ColumnWidth = 13.6666666666667
$ws.Columns("L").ColumnWidth = 15.8333333333333
$ws.Columns("M:N").ColumnWidth = 16.3333333333333
$ws.Columns("O").ColumnWidth = 16
$ws.Columns("P").ColumnWidth = 16.3333333333333
$ws.Columns("Q").ColumnWidth = 16
$ws.Columns("R").ColumnWidth = 11.8333333333333
$ws.Columns("S").ColumnWidth = 20.5
$ws.Columns("T").ColumnWidth = 15
$ws.Columns("U").ColumnWidth = 20.5
$ws.Columns("V").ColumnWidth = 16
$ws.Columns("W:X").ColumnWidth = 15
$ws.Columns("Y").ColumnWidth = 16
